$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9998340106875503
$ws.Range("C2").Value = 0.2532704065238534
$ws.Range("D2").Value = 0.02824048846233751
$ws.Range("E2").Value = 0.1184694903238031
$ws.Range("F2").Value = 0.8287855757083662
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("L2").Value = 0.1886793512131248
$ws.Range("M2").Value = 0.2225080754306248
$ws.Range("N2").Value = 1.376402873543597
$ws.Range("O2").Value = 2.873585931597034
$ws.Range("B3").Value = 0.9154134562347451
$ws.Range("C3").Value = 0.2429535528182356
$ws.Range("D3").Value = 0.0270744962661027
$ws.Range("E3").Value = 0.1194771694374737
$ws.Range("F3").Value = 0.8230670106637845
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("L3").Value = 0.1860235826452055
$ws.Range("M3").Value = 0.2092003295402449
$ws.Range("N3").Value = 1.390789325767418
$ws.Range("O3").Value = 2.869311570831997
$ws.Range("B4").Value = 0.8637903451062812
$ws.Range("C4").Value = 0.2365655221593812
$ws.Range("D4").Value = 0.02635243062252002
$ws.Range("E4").Value = 0.1201334435495065
$ws.Range("F4").Value = 0.8201120332256338
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("L4").Value = 0.1844843663127094
$ws.Range("M4").Value = 0.2011066445347964
$ws.Range("N4").Value = 1.400155429781762
$ws.Range("O4").Value = 2.868562608234441
$ws.Range("B5").Value = 0.8428078309324576
$ws.Range("C5").Value = 0.2339490491903433
$ws.Range("D5").Value = 0.02605665567152116
$ws.Range("E5").Value = 0.1204103349571133
$ws.Range("F5").Value = 0.8190476355548668
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("L5").Value = 0.183880166847679
$ws.Range("M5").Value = 0.1978280302640343
$ws.Range("N5").Value = 1.404106192046012
$ws.Range("O5").Value = 2.8687286328487
$ws.Range("B6").Value = 0.8393270153740389
$ws.Range("C6").Value = 0.2335137867586354
$ws.Range("D6").Value = 0.02600745072161459
$ws.Range("E6").Value = 0.1204568839050189
$ws.Range("F6").Value = 0.8188793327823376
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("L6").Value = 0.18378123328295
$ws.Range("M6").Value = 0.1972848097824738
$ws.Range("N6").Value = 1.404770307293489
$ws.Range("O6").Value = 2.868784651350154
$ws.Range("B7").Value = 0.8635071461037853
$ws.Range("C7").Value = 0.2365302891501102
$ws.Range("D7").Value = 0.0263484478590108
$ws.Range("E7").Value = 0.1201371395088011
$ws.Range("F7").Value = 0.8200971125481047
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("L7").Value = 0.1844761244960225
$ws.Range("M7").Value = 0.2010623482783629
$ws.Range("N7").Value = 1.400208168545515
$ws.Range("O7").Value = 2.8685629398054
$ws.Range("B8").Value = 0.9706825803170887
$ws.Range("C8").Value = 0.2497243277062609
$ws.Range("D8").Value = 0.02783973895927971
$ws.Range("E8").Value = 0.1188091501259567
$ws.Range("F8").Value = 0.8266983059276853
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("L8").Value = 0.1877446940578054
$ws.Range("M8").Value = 0.217903621792324
$ws.Range("N8").Value = 1.381252765874365
$ws.Range("O8").Value = 2.871722630734695
$ws.Range("B9").Value = 1.182491285369963
$ws.Range("C9").Value = 0.2751687732997254
$ws.Range("D9").Value = 0.03071480266731896
$ws.Range("E9").Value = 0.1165024395906844
$ws.Range("F9").Value = 0.844062936668692
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("L9").Value = 0.1948781482326112
$ws.Range("M9").Value = 0.2515365702152224
$ws.Range("N9").Value = 1.348306900267936
$ws.Range("O9").Value = 2.89282208470641
$ws.Range("B10").Value = 1.339068263590491
$ws.Range("C10").Value = 0.2935962198798734
$ws.Range("D10").Value = 0.03279637627263554
$ws.Range("E10").Value = 0.1149883120733852
$ws.Range("F10").Value = 0.8595264885623237
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("L10").Value = 0.2005588575999298
$ws.Range("M10").Value = 0.2766111450380748
$ws.Range("N10").Value = 1.326675052089328
$ws.Range("O10").Value = 2.91744667579394
$ws.Range("B11").Value = 1.410500935366201
$ws.Range("C11").Value = 0.3019204887245053
$ws.Range("D11").Value = 0.03373653787132014
$ws.Range("E11").Value = 0.1143385579397709
$ws.Range("F11").Value = 0.8671514069051796
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("L11").Value = 0.2032383865872021
$ws.Range("M11").Value = 0.2880963012296149
$ws.Range("N11").Value = 1.317392492080685
$ws.Range("O11").Value = 2.930638455276153
$ws.Range("B12").Value = 1.437579164946953
$ws.Range("C12").Value = 0.3050641451535796
$ws.Range("D12").Value = 0.0340915660672394
$ws.Range("E12").Value = 0.1140981144268842
$ws.Range("F12").Value = 0.8701238267605333
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("L12").Value = 0.2042667260072903
$ws.Range("M12").Value = 0.2924565858272317
$ws.Range("N12").Value = 1.313957642468303
$ws.Range("O12").Value = 2.935920541193411
$ws.Range("B13").Value = 1.431746148260913
$ws.Range("C13").Value = 0.3043874858347237
$ws.Range("D13").Value = 0.03401514875208989
$ws.Range("E13").Value = 0.1141496491176257
$ws.Range("F13").Value = 0.8694798793890186
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("L13").Value = 0.2040446479907558
$ws.Range("M13").Value = 0.2915170297712208
$ws.Range("N13").Value = 1.314693828891144
$ws.Range("O13").Value = 2.934770192977794
$ws.Range("B14").Value = 1.412728119076803
$ws.Range("C14").Value = 0.3021792917466826
$ws.Range("D14").Value = 0.03376576621617033
$ws.Range("E14").Value = 0.1143186642358724
$ws.Range("F14").Value = 0.8673942448907042
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("L14").Value = 0.2033227152141137
$ws.Range("M14").Value = 0.2884548030377161
$ws.Range("N14").Value = 1.317108296317201
$ws.Range("O14").Value = 2.931067268030915
$ws.Range("B15").Value = 1.401082665866568
$ws.Range("C15").Value = 0.3008255892684986
$ws.Range("D15").Value = 0.0336128826439861
$ws.Range("E15").Value = 0.1144229204761213
$ws.Range("F15").Value = 0.8661278103284218
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("L15").Value = 0.2028822879686913
$ws.Range("M15").Value = 0.2865805410885116
$ws.Range("N15").Value = 1.318597679860282
$ws.Range("O15").Value = 2.928836463947761
$ws.Range("B16").Value = 1.334403988728411
$ws.Range("C16").Value = 0.2930510194657643
$ws.Range("D16").Value = 0.03273479691502956
$ws.Range("E16").Value = 0.1150315597161082
$ws.Range("F16").Value = 0.8590400742649393
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("L16").Value = 0.2003856583486652
$ws.Range("M16").Value = 0.2758621278752145
$ws.Range("N16").Value = 1.327292917792448
$ws.Range("O16").Value = 2.916624646883633
$ws.Range("B17").Value = 1.293550397174613
$ws.Range("C17").Value = 0.2882664864481512
$ws.Range("D17").Value = 0.03219437568802874
$ws.Range("E17").Value = 0.1154149308537473
$ws.Range("F17").Value = 0.8548432971164885
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("L17").Value = 0.1988784412296809
$ws.Range("M17").Value = 0.2693067203457957
$ws.Range("N17").Value = 1.332770071880852
$ws.Range("O17").Value = 2.909643114798598
$ws.Range("B18").Value = 1.270071864073543
$ws.Range("C18").Value = 0.2855090555095501
$ws.Range("D18").Value = 0.03188290485441314
$ws.Range("E18").Value = 0.1156391099730243
$ws.Range("F18").Value = 0.8524849927366915
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("L18").Value = 0.1980205072332382
$ws.Range("M18").Value = 0.2655436339572717
$ws.Range("N18").Value = 1.335972893777416
$ws.Range("O18").Value = 2.905814794237386
$ws.Range("B19").Value = 1.262125801635136
$ws.Range("C19").Value = 0.2845744979084373
$ws.Range("D19").Value = 0.03177733772768221
$ws.Range("E19").Value = 0.1157156444659355
$ws.Range("F19").Value = 0.8516960520747858
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("L19").Value = 0.1977315689491803
$ws.Range("M19").Value = 0.2642707963950954
$ws.Range("N19").Value = 1.337066331576054
$ws.Range("O19").Value = 2.904550739695821
$ws.Range("B20").Value = 1.297897335878588
$ws.Range("C20").Value = 0.2887763779036447
$ws.Range("D20").Value = 0.03225197029090054
$ws.Range("E20").Value = 0.1153737401375965
$ws.Range("F20").Value = 0.8552842989634684
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("L20").Value = 0.1990379583026254
$ws.Range("M20").Value = 0.270003789132204
$ws.Range("N20").Value = 1.332181584900454
$ws.Range("O20").Value = 2.910366925250401
$ws.Range("B21").Value = 1.418313419165031
$ws.Range("C21").Value = 0.3028281255286629
$ws.Range("D21").Value = 0.03383904294785367
$ws.Range("E21").Value = 0.1142688683651452
$ws.Range("F21").Value = 0.8680045376944463
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("L21").Value = 0.2035343940475229
$ws.Range("M21").Value = 0.2893539534457119
$ws.Range("N21").Value = 1.316396930345149
$ws.Range("O21").Value = 2.932147123433793
$ws.Range("B22").Value = 1.497176238783084
$ws.Range("C22").Value = 0.3119617681342959
$ws.Range("D22").Value = 0.03487050392620716
$ws.Range("E22").Value = 0.1135794330497442
$ws.Range("F22").Value = 0.8768136165093097
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("L22").Value = 0.2065526718712647
$ws.Range("M22").Value = 0.3020650422218765
$ws.Range("N22").Value = 1.306548548254597
$ws.Range("O22").Value = 2.948052759101955
$ws.Range("B23").Value = 1.455071080470532
$ws.Range("C23").Value = 0.3070915963349421
$ws.Range("D23").Value = 0.03432052917890616
$ws.Range("E23").Value = 0.1139444117435096
$ws.Range("F23").Value = 0.8720666524777698
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("L23").Value = 0.2049344937795468
$ws.Range("M23").Value = 0.2952750460107012
$ws.Range("N23").Value = 1.311761999376245
$ws.Range("O23").Value = 2.9394105604907
$ws.Range("B24").Value = 1.295932057568905
$ws.Range("C24").Value = 0.2885458769584091
$ws.Range("D24").Value = 0.03222593418489339
$ws.Range("E24").Value = 0.115392350704357
$ws.Range("F24").Value = 0.8550847523389677
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("L24").Value = 0.1989658138886767
$ws.Range("M24").Value = 0.269688626549538
$ws.Range("N24").Value = 1.332447471864668
$ws.Range("O24").Value = 2.910039112894935
$ws.Range("B25").Value = 1.125019938851722
$ws.Range("C25").Value = 0.2683318293564412
$ws.Range("D25").Value = 0.02994237257370003
$ws.Range("E25").Value = 0.1170946956531096
$ws.Range("F25").Value = 0.8388910036447186
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("L25").Value = 0.1928710261636155
$ws.Range("M25").Value = 0.242373551117403
$ws.Range("N25").Value = 1.356767530887481
$ws.Range("O25").Value = 2.885514754911753
